$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the last data row (111) -------------------------------------
# The refreshed stock export has one fewer SKU row than before; removing it
# also lets Excel reclaim the now-unused "Линия Консумаш" shared string.
$ws.Rows.Item(111).Delete()

# --- 2. Refresh every remaining data row (2-110) with the new figures ----
$data = @(
    @('ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы', 65888),
    @('Сб. Фитонефрол (Урологический сбор) 50г', 10139),
    @('Рябина плоды 50г', 252),
    @('Дуба кора 75г', 49638),
    @('Полынь горькая трава 50г', 45468),
    @('Аир корневища 75г', 7681),
    @('Укроп пахучий плоды 50г', 72766),
    @('Череда трава 50г', 12896),
    @('Мята перечная листья 50г', 16436),
    @('Крапива листья 50г', 15128),
    @('Бессмертник песчаный цветки 30г', 31497),
    @('Эрва шерстистая трава 30г', 20635),
    @('Чага (березовый гриб) 50г', 29616),
    @('Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г', 6587),
    @('Лен семена 100г', 97580),
    @('Ромашка цветки вн 50г', 102002),
    @('Кукуруза столбики с рыльцами 40г', 27615),
    @('Ламинарии слоевища (морская капуста) 100г', 17245),
    @('Можжевельник плоды 50г', 13198),
    @('Солодка корни 50г', 25726),
    @('Пустырник трава 50г', 13219),
    @('Сенна листья 50г', 34467),
    @('Брусника листья 50г', 19225),
    @('Крушина кора 50г', 15668),
    @('Девясил корневища и корни 50г', 15554),
    @('Спорыш трава 50г', 19391),
    @('Зверобой трава 50г', 47331),
    @('Валериана корневища с корнями 50г', 22948),
    @('Чабрец трава 50г', 25211),
    @('Сб. Грудной №4 50г', 41949),
    @('Сб. Фитопектол №1 (Грудной сбор №1) 35г', 5576),
    @('Сб. Фитопектол №2 (Грудной сбор №2) 35г', 9085),
    @('Шиповник плоды низковитаминные 50г', 34826),
    @('Чистотел трава 50г', 29927),
    @('Эвкалипт прутовидный листья 75г', 23604),
    @('Ноготки цветки 50г', 30782),
    @('Тысячелистник трава 50г', 23449),
    @('Толокнянка листья 50г', 11434),
    @('Пижма цветки 75г', 27653),
    @('Мать-и-мачеха листья 35г', 31979),
    @('Багульник болотный побеги 50г', 19516),
    @('Боярышник плоды 75г', 32199),
    @('Шалфей листья 50г', 45834),
    @('Подорожник большой листья 50г', 16810),
    @('Алтей корни 75г', 8075),
    @('Береза почки 50г', 24668),
    @('Липа цветки 35г', 25631),
    @('Фп "ФармаЦветик® Фиточай для кормящих мам" 20х1,5 г', 2420),
    @('Фп Фиточай "Лактафитол" (БАД) 20х1,5 г', 12506),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г', 2340),
    @('Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г', 4100),
    @('Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г', 2050),
    @('Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г', 2460),
    @('Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г', 3820),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г', 5000),
    @('Фп Пижма цветки 20х1,5г', 2442),
    @('Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г', 107520),
    @('Фп Сб. Бруснивер 20x2,0г', 145513),
    @('Фп Зверобой трава 20x1,5г', 28092),
    @('Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г', 41991),
    @('Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г', 56589),
    @('Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г', 4464),
    @('Фп Мелисса лекарственная трава 20x1,5г', 28136),
    @('Фп Сенна листья 20x1,5г', 80238),
    @('Фп Мята перечная листья 20x1,5г', 39987),
    @('Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г', 53246),
    @('Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г', 5436),
    @('Фп Сб. Желудочный №3 20x2,0г', 14958),
    @('Фп Толокнянка листья 20x1,5г', 27647),
    @('Фп Хвощ полевой трава 20х1,5г', 26158),
    @('Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г', 4194),
    @('Фп Брусника листья 20х1,5г', 60313),
    @('Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г', 18126),
    @('Фп Пастушья сумка трава 20х1,5г', 5398),
    @('Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г', 4446),
    @('Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г', 3564),
    @('Фп Ромашка цветки 20x1,5г', 975634),
    @('Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г', 42047),
    @('Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г', 4554),
    @('Фп Сб. Арфазетин-Э 20x2,0г', 49366),
    @('Фп Шалфей листья 20х1,5г', 84177),
    @('Фп Череда трава 20х1,5г', 63369),
    @('Фп Крапива листья 20x1,5г', 60298),
    @('Фп Шиповник плоды 20х2,0г', 38568),
    @('Фп Подорожник листья 20x1,5г', 30159),
    @('Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г', 2070),
    @('Фп Чабрец трава 20x1,5 г', 58672),
    @('Фп Чистотел трава 20х1,5г', 38709),
    @('Фп Сб. Элекасол 20x2,0г', 31590),
    @('Фп Ольха соплодия 20х1,5г', 3308),
    @('Фп Береза листья 20x1,5г', 5163),
    @('Фп Липа цветки 20x1,5г', 54660),
    @('Фп Душица трава 20x1,5г', 32764),
    @('Фп Пустырник трава 20x1,5г', 41885),
    @('Фп Золототысячник трава 20х1,5г', 4726),
    @('Фп Боярышник плоды 20х3,0г', 17496),
    @('Фп Сб. Грудной №4 20x2,0г', 834318),
    @('Фп Фиалка трехцветная трава 20x1,5г', 7182),
    @('Фп Аир корневища 20x1,5г', 8662),
    @('Фп Дуб кора 20х1,5г', 6183),
    @('Фп Почечный чай листья 20x1,5г', 63165),
    @('Фп Лапчатка корневища 20x2,5г', 3290),
    @('Фп Кровохлебка корневища и корни 20x1,5г', 7858),
    @('Фп Тысячелистник трава 20x1,5г', 25495),
    @('Фп Валериана корневища с корнями 20x1,5г', 22123),
    @('Фп Ноготки цветки 20x1,5г', 61574),
    @('Фп Крушина кора 20x1,5г', 21864),
    @('Фп Девясил корневища и корни 20х1,5г', 21669),
    @('Фп Бадан корневища 20x1,5г', 5291)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# --- 3. Row 4 now carries the one-off integer (no thousands separator) --
#     number format that used to live on row 58; move the live format
#     there first so no brand-new style slot gets allocated ...
$ws.Range("B58").Copy()
$ws.Range("B4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ... then restore row 58 (and the former bold separator row, now a plain
#     data row) to the standard look used by every other quantity cell.
$ws.Range("B2").Copy()
$ws.Range("B58").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2:B2").Copy()
$ws.Range("A101:B101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- 4. Restore the sheet dimension / view -------------------------------
$ws.Range("A94").Select()
